$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Ordered (old, new) pairs for each cell, walked in row-major order
# (matches the natural Tables.Item(1).Cell(r,c) traversal).
$replacements = @(
    ,@("89-67=", "78+17=")
    ,@("20+56=", "7+6=")
    ,@("18+38=", "74-28=")
    ,@("31+20=", "22+63=")
    ,@("24+70=", "69-18=")
    ,@("23+68=", "11+74=")
    ,@("71-12=", "19+45=")
    ,@("71-31=", "67-28=")
    ,@("39+56=", "87-0=")
    ,@("6+60=", "39+36=")
    ,@("31-22=", "40-22=")
    ,@("86-48=", "80-3=")
    ,@("46+23=", "21+27=")
    ,@("6+60=", "65-28=")
    ,@("41+42=", "33+41=")
    ,@("26+32=", "31+30=")
    ,@("21+19=", "95-93=")
    ,@("16-2=", "38-12=")
    ,@("7+66=", "2+2=")
    ,@("41+10=", "42+11=")
    ,@("52-44=", "10+77=")
    ,@("74-3=", "98-5=")
    ,@("87+4=", "44-21=")
    ,@("35+55=", "81+1=")
    ,@("30+26=", "59-19=")
    ,@("38-14=", "38-19=")
    ,@("38+43=", "40+51=")
    ,@("18+35=", "66-41=")
    ,@("99-72=", "79+16=")
    ,@("67-39=", "97-21=")
    ,@("47-19=", "80+11=")
    ,@("94-79=", "9+69=")
    ,@("8+48=", "58+32=")
    ,@("59-28=", "15-8=")
    ,@("96-58=", "91-87=")
    ,@("85-62=", "21+77=")
    ,@("51-40=", "16+66=")
    ,@("99-14=", "47-6=")
    ,@("64-1=", "61+37=")
    ,@("85-62=", "74-24=")
    ,@("48+33=", "62+18=")
    ,@("78-33=", "59+26=")
    ,@("54+9=", "57+19=")
    ,@("78-30=", "90-66=")
    ,@("25+16=", "22+39=")
    ,@("1+45=", "32+63=")
    ,@("75-17=", "59+23=")
    ,@("97-8=", "3+96=")
    ,@("91-90=", "96-34=")
    ,@("96-53=", "38-32=")
    ,@("86+9=", "44+36=")
    ,@("41+8=", "76-24=")
    ,@("77-32=", "33+55=")
    ,@("21+23=", "65+31=")
    ,@("34+48=", "94-34=")
    ,@("43+27=", "48+12=")
    ,@("70+8=", "11+25=")
    ,@("83+14=", "58-45=")
    ,@("1+7=", "43+22=")
    ,@("37-14=", "60-15=")
    ,@("46+16=", "80-10=")
    ,@("3+47=", "68+31=")
    ,@("56-42=", "68+31=")
    ,@("95-92=", "0+17=")
    ,@("55-37=", "0+9=")
    ,@("5+16=", "9+86=")
    ,@("66-38=", "17-6=")
    ,@("18+44=", "48+16=")
    ,@("82+9=", "75-30=")
    ,@("83-52=", "41+9=")
    ,@("53-44=", "56+11=")
    ,@("81-69=", "94-76=")
    ,@("95-24=", "26+14=")
    ,@("4+37=", "27+26=")
    ,@("66-56=", "23+38=")
    ,@("79-55=", "2+79=")
    ,@("49+48=", "39-27=")
    ,@("53-47=", "17+70=")
    ,@("30+59=", "94-21=")
    ,@("84-38=", "3+20=")
    ,@("70-53=", "33+40=")
    ,@("97-62=", "67+19=")
    ,@("88-58=", "84-14=")
    ,@("20+33=", "67+22=")
    ,@("27+48=", "12+79=")
    ,@("27+2=", "42-27=")
    ,@("47-22=", "78-12=")
    ,@("37+28=", "82-68=")
    ,@("6+20=", "36+62=")
    ,@("16+44=", "92-16=")
    ,@("16+82=", "76-9=")
    ,@("18+2=", "60+8=")
    ,@("22+69=", "55-12=")
    ,@("32+62=", "35+17=")
    ,@("53-20=", "76-41=")
    ,@("34+23=", "83-71=")
    ,@("13+69=", "87-34=")
    ,@("44+52=", "50+12=")
    ,@("28+6=", "22+47=")
    ,@("99-2=", "98-83=")
)

$rows = $t.Rows.Count
$cols = $t.Columns.Count
$i = 0
$mismatches = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $old = $replacements[$i][0]
        $new = $replacements[$i][1]
        $cell = $t.Cell($r, $c)
        $cur = $cell.Range.Text
        $curTrim = $cur.Substring(0, $cur.Length - 2)
        if ($curTrim -ne $old) {
            Write-Host "MISMATCH at ($r,$c): expected [$old] got [$curTrim]"
            $mismatches = $mismatches + 1
        }
        $cell.Range.Text = $new
        $i = $i + 1
    }
}
Write-Host "Processed $i cells, $mismatches mismatches"
